# "video 1 parte 2"
# Insert two new columns (A, B) at the left of the sheet, add a "fila" /
# "resultado" status pair for each data row, and keep the existing
# hyperlinks (on the "correo" column, now column G) pointing at their
# original mailto: targets.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing columns A:J two to the right -> C:L
$ws.Columns("A:B").Insert()

# Column-insert doesn't relocate the worksheet's hyperlink anchors in this
# engine, so rebuild them explicitly at their new location (G2:G7),
# preserving the original rIds / targets (mailto:jeisson@gmail.com /
# mailto:luna@gmail.com alternating, matching the pre-edit E2:E7 links).
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("G2"), "mailto:jeisson@gmail.com")
$ws.Hyperlinks.Add($ws.Range("G3"), "mailto:luna@gmail.com")
$ws.Hyperlinks.Add($ws.Range("G4"), "mailto:jeisson@gmail.com")
$ws.Hyperlinks.Add($ws.Range("G5"), "mailto:luna@gmail.com")
$ws.Hyperlinks.Add($ws.Range("G6"), "mailto:jeisson@gmail.com")
$ws.Hyperlinks.Add($ws.Range("G7"), "mailto:luna@gmail.com")
# Adding a hyperlink re-styles the cell with a fresh "Hipervinculo" xf;
# restore the original text number format so the cell collapses back onto
# the same style the sheet already used for these mail cells.
$ws.Range("G2:G7").NumberFormat = "@"

# New leading columns: "fila" (row number) and "resultado" ("ok" for
# every data row).
$ws.Range("A1:B7").Style = "Normal"
$ws.Range("A1").Value = "fila"
$ws.Range("B1").Value = "resultado"
for ($r = 2; $r -le 7; $r++) {
    $ws.Cells.Item($r, 1).Value = $r - 1
    $ws.Cells.Item($r, 2).Value = "ok"
}

# Match the author's final selection in the saved file.
$null = $ws.Range("L5").Select()
